$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "60.005.71"
Set-TextValue $ws.Range("E2") "  -0.64%  "

Set-TextValue $ws.Range("D3") "2.418.21"
Set-TextValue $ws.Range("E3") "  -1.05%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "553.12"
Set-TextValue $ws.Range("E5") "  -0.49%  "

Set-TextValue $ws.Range("D6") "137.12"
Set-TextValue $ws.Range("E6") "  -1.07%  "

Set-TextValue $ws.Range("D8") "0.598"
Set-TextValue $ws.Range("E8") "  +4.78%  "

Set-TextValue $ws.Range("E9") "  -1.19%  "

Set-TextValue $ws.Range("D10") "5.67"
Set-TextValue $ws.Range("E10") "  -1.92%  "

Set-TextValue $ws.Range("E12") "  -1.80%  "

Set-TextValue $ws.Range("D13") "25.27"
Set-TextValue $ws.Range("E13") "  +1.51%  "

Set-TextValue $ws.Range("D14") "2.848.80"
Set-TextValue $ws.Range("E14") "  -1.04%  "

Set-TextValue $ws.Range("D15") "59.944.93"
Set-TextValue $ws.Range("E15") "  -0.62%  "

Set-TextValue $ws.Range("D16") "0.0000137"
Set-TextValue $ws.Range("E16") "  -1.59%  "

Set-TextValue $ws.Range("D17") "2.413.89"
Set-TextValue $ws.Range("E17") "  -1.86%  "

Set-TextValue $ws.Range("D18") "11.28"
Set-TextValue $ws.Range("E18") "  -1.01%  "

Set-TextValue $ws.Range("D19") "4.44"
Set-TextValue $ws.Range("E19") "  +0.41%  "

Set-TextValue $ws.Range("D20") "327.94"
Set-TextValue $ws.Range("E20") "  -2.21%  "

Set-TextValue $ws.Range("E21") "  -3.67%  "

Set-TextValue $ws.Range("E22") "  +0.07%  "

Set-TextValue $ws.Range("D23") "66.02"
Set-TextValue $ws.Range("E23") "  +2.08%  "

Set-TextValue $ws.Range("D24") "0.177"
Set-TextValue $ws.Range("E24") "  +3.34%  "

Set-TextValue $ws.Range("E25") "  +0.45%  "

Set-TextValue $ws.Range("E26") "  +0.07%  "

Set-TextValue $ws.Range("D28") "0.0₃0775"
Set-TextValue $ws.Range("E28") "  -2.02%  "

Set-TextValue $ws.Range("E29") "  -1.61%  "

Set-TextValue $ws.Range("D30") "168.39"
Set-TextValue $ws.Range("E30") "  -1.69%  "

Set-TextValue $ws.Range("D31") "6.06"
Set-TextValue $ws.Range("E31") "  -4.02%  "

Set-TextValue $ws.Range("E32") "  +1.49%  "

Set-TextValue $ws.Range("D33") "18.60"
Set-TextValue $ws.Range("E33") "  -1.26%  "

Set-TextValue $ws.Range("E35") "  -0.54%  "

Set-TextValue $ws.Range("E36") "  +0.01%  "

Set-TextValue $ws.Range("D37") "4.19"
Set-TextValue $ws.Range("E37") "  -1.87%  "

Set-TextValue $ws.Range("D38") "325.57"
Set-TextValue $ws.Range("E38") "  +2.37%  "

Set-TextValue $ws.Range("E39") "  -2.08%  "

Set-TextValue $ws.Range("D40") "0.405"
Set-TextValue $ws.Range("E40") "  -2.21%  "

Set-TextValue $ws.Range("E41") "  -1.79%  "

Set-TextValue $ws.Range("D42") "140.55"
Set-TextValue $ws.Range("E42") "  -2.43%  "

Set-TextValue $ws.Range("E43") "  +0.91%  "

Set-TextValue $ws.Range("D44") "19.61"
Set-TextValue $ws.Range("E44") "  -1.70%  "

Set-TextValue $ws.Range("D45") "0.0518"
Set-TextValue $ws.Range("E45") "  -1.31%  "

Set-TextValue $ws.Range("D46") "0.576"
Set-TextValue $ws.Range("E46") "  +0.47%  "

Set-TextValue $ws.Range("D47") "0.398"
Set-TextValue $ws.Range("E47") "  -2.89%  "

Set-TextValue $ws.Range("E49") "  -0.08%  "

Set-TextValue $ws.Range("E50") "  -4.09%  "

Set-TextValue $ws.Range("D51") "4.66"
Set-TextValue $ws.Range("E51") "  -1.12%  "
